$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.29429999999999
$ws.Range("A9").Value = -21.6673
$ws.Range("C9").Value = -10.46
$ws.Range("A18").Value = -22.26390000000001
$ws.Range("A20").Value = -20.28609999999999
$ws.Range("C23").Value = -12.34310000000001
$ws.Range("C24").Value = -12.8048
$ws.Range("C26").Value = -12.34050000000001
$ws.Range("A27").Value = -22.00680000000001
$ws.Range("C34").Value = -11.79960000000001
$ws.Range("C35").Value = -12.58810000000001
$ws.Range("C48").Value = -12.1442
$ws.Range("C52").Value = -11.2142
$ws.Range("C66").Value = -11.6116
$ws.Range("C67").Value = -11.5533
$ws.Range("A69").Value = -21.47939999999997
$ws.Range("A76").Value = -19.68449999999999
$ws.Range("C80").Value = -13.13190000000001
$ws.Range("A82").Value = -21.7859
$ws.Range("C99").Value = -12.58829999999999
